$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph whose text (without trailing paragraph mark)
# equals the given string.
# ---------------------------------------------------------------------------
function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# Helper: insert a new italic paragraph right after the given paragraph.
function Insert-ItalicParagraphAfter($paragraph, $newText) {
    $paragraph.Range.InsertParagraphAfter()
    $newRange = $paragraph.Next().Range
    $newRange.MoveEnd(1, -1)
    $newRange.Text = $newText
    $newRange.Font.Italic = $true
}

# ---------------------------------------------------------------------------
# 1. Update the activation date.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# ---------------------------------------------------------------------------
# 2. Insert the English translation of "Objetivos".
# ---------------------------------------------------------------------------
$p = Find-ParagraphByText $d "A disciplina visa propiciar aos alunos os conhecimentos básicos de eletroquímica, tanto do ponto de vista da eletroquímica iônica como da eletródica, e apresentar as principais aplicações da eletroquímica"
Insert-ItalicParagraphAfter $p "The course aims to provide students with basic knowledge of electrochemistry, both from the point of view of ionic and electrodic electrochemistry, and to present the main applications of electrochemistry"

# ---------------------------------------------------------------------------
# 3. Insert the English translation of "Programa resumido".
# ---------------------------------------------------------------------------
$p = Find-ParagraphByText $d "Princípios da eletroquímica iônica e da eletroquímica eletródica. Aplicações."
Insert-ItalicParagraphAfter $p "Principles of ionic electrochemistry and electrodic electrochemistry. Applications."

# ---------------------------------------------------------------------------
# 4. Insert the English translation of "Programa".
# ---------------------------------------------------------------------------
$p = Find-ParagraphByText $d "Princípios da eletroquímica iônica: interações iônicas, equilíbrio iônico e condução eletrolítica. Princípios da eletroquímica eletródica: fenômenos interfaciais, potenciais de eletrodo e células eletroquímicas. Processos de eletrodo. Métodos eletroquímicos de análise química. Aplicações da eletroquímica: fontes eletroquímicas de energia, processos eletrometalúrgicos e galvanoplastia."
Insert-ItalicParagraphAfter $p "Principles of ionic electrochemistry: ionic interactions, ionic equilibrium and electrolyte conduction. Principles of electrochemical electrochemistry: interfacial phenomena, electrode potentials and electrochemical cells. Electrode processes. Electrochemical methods of chemical analysis. Applications of electrochemistry: electrochemical sources of energy, electrometallurgical processes and electroplating."

Write-Host "Edits applied successfully."
